$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting existing rows 47:173 down to 48:174.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new record's data.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R carry over the same constant
# values used throughout this data block; only D, J, K, L, M, P are new.
$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value = "La Araucanía"
$ws.Cells.Item(47, 4).Value = 44979
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(47, 6).Value = 100114002
$ws.Cells.Item(47, 7).Value = "Camote"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 5
$ws.Cells.Item(47, 11).Value = 26000
$ws.Cells.Item(47, 12).Value = 26000
$ws.Cells.Item(47, 13).Value = 26000
$ws.Cells.Item(47, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(47, 15).Value = "Perú"
$ws.Cells.Item(47, 16).Value = 1300
$ws.Cells.Item(47, 17).Value = 20
$ws.Cells.Item(47, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(47, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
